$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.103.47'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.98%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.243.43'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.98%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.57'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.81%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.633'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '77.04'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.62%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.632'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.09%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.14'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.25%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0958'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.29%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.22'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.66%  '
$ws.Range("E13").Value = '  -2.86%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.577.99'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.11%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.88'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.27%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.862'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.19%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.256.59'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.87%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '41.975.18'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.07%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0984'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.98%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.13'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.89%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.86'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.83%  '
$ws.Range("E22").Value = '  +4.83%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '231.82'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.35%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.43'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.76%  '
$ws.Range("E25").Value = '  +0.09%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.68'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.73%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.31'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.52%  '
$ws.Range("E28").Value = '  +12.58%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.16'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.50%  '
$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '168.95'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.18%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.59'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.29%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '33.15'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.55%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0833'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.19%  '
$ws.Range("E34").Value = '  -5.52%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.125'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.28%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.53'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.01%  '
$ws.Range("E37").Value = '  +2.87%  '
$ws.Range("E38").Value = '  -2.46%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '14.21'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.42%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.93'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.03%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.20'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.32%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '113.27'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +14.12%  '
$ws.Range("E43").Value = '  -6.75%  '
$ws.Range("E44").Value = '  -0.72%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.72'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.04%  '
$ws.Range("E46").Value = '  -2.89%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.997'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.40%  '
$ws.Range("E48").Value = '  -2.95%  '
$ws.Range("E49").Value = '  -1.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.26'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -13.36%  '
$ws.Range("E51").Value = '  -1.91%  '
